$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Two brand-new rows ("line7", "line8") are inserted into the table right
# after "line6", pushing the "extr" rows down by two positions; the two
# displaced rows ("extr7", "extr8") end up in two newly appended rows
# (16 and 17).
#
# To make sure the workbook's shared-string table ends up in the exact same
# (fresh, sequential) order as the target file - name, from_bus, to_bus,
# in_service, line1..line8, extr1..extr8 - we rebuild all the cell VALUES
# from scratch, in that exact left-to-right / top-to-bottom order. Plain
# in-place edits would instead leave "line7"/"line8" appended at the *end*
# of the shared-string table (since Excel only ever appends brand new
# strings), so a full, ordered re-entry is required to reproduce the table
# layout exactly.

# 1) Extend the existing formatting (bold/border/centered style used by the
#    "name" column) down to the two new rows by copying the format of the
#    last existing row - this does not disturb the style table.
$ws.Range("A15:E15").Copy($ws.Range("A16:E16"))
$ws.Range("A15:E15").Copy($ws.Range("A17:E17"))

# 2) Wipe all values (but keep formatting) so the shared-string table is
#    rebuilt from scratch as we re-type the values below.
$ws.Range("A1:E17").ClearContents()

# 3) Re-enter every value, row by row, in the final desired layout.
$ws.Range("B1").Value = "name"
$ws.Range("C1").Value = "from_bus"
$ws.Range("D1").Value = "to_bus"
$ws.Range("E1").Value = "in_service"

$data = @(
    @(0,  "line1", 7,  9,  $false),
    @(1,  "line2", 9,  8,  $true),
    @(2,  "line3", 8,  10, $true),
    @(3,  "line4", 8,  11, $true),
    @(4,  "line5", 10, 5,  $true),
    @(5,  "line6", 12, 8,  $true),
    @(6,  "line7", 14, 11, $true),
    @(7,  "line8", 16, 9,  $true),
    @(8,  "extr1", 5,  12, $true),
    @(9,  "extr2", 5,  9,  $true),
    @(10, "extr3", 10, 11, $true),
    @(11, "extr4", 7,  8,  $false),
    @(12, "extr5", 9,  11, $true),
    @(13, "extr6", 7,  11, $false),
    @(14, "extr7", 5,  7,  $true),
    @(15, "extr8", 8,  5,  $false)
)

$row = 2
foreach ($r in $data) {
    $ws.Cells.Item($row, 1).Value = $r[0]
    $ws.Cells.Item($row, 2).Value = $r[1]
    $ws.Cells.Item($row, 3).Value = $r[2]
    $ws.Cells.Item($row, 4).Value = $r[3]
    $ws.Cells.Item($row, 5).Value = $r[4]
    $row++
}
